$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Some Price values are plain decimal numbers stored as TEXT in the sheet
# (e.g. "219.33"); force those cells to text format first so Excel does not
# reinterpret them as numbers when the new value is assigned.
$ws.Range("D2").Value = "28.307.06"
$ws.Range("E2").Value = "  +3.95%  "
$ws.Range("D3").Value = "1.731.01"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.33"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.524"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.09"
$ws.Range("E8").Value = "  +3.93%  "
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0895"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "1.976.01"
$ws.Range("E12").Value = "  +2.49%  "
$ws.Range("D13").Value = "1.730.95"
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.25"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.566"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.77"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "28.298.97"
$ws.Range("E17").Value = "  +3.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "245.91"
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +1.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.70"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.05"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.49"
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.68"
$ws.Range("E27").Value = "  +1.27%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0518"
$ws.Range("E30").Value = "  +2.87%  "
$ws.Range("E31").Value = "  +2.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.42"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").Value = "1.484.22"
$ws.Range("E34").Value = "  -3.99%  "
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.986"
$ws.Range("E36").Value = "  +3.71%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.92"
$ws.Range("E41").Value = "  +0.85%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").Value = "1.880.52"
$ws.Range("E44").Value = "  +2.35%  "
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.802"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.73"
$ws.Range("E47").Value = "  +7.57%  "
$ws.Range("D48").Value = "0.0₆0114"
$ws.Range("E48").Value = "  +3.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "90.35"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.20"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("E51").Value = "  -0.67%  "
